$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the date number format used by the existing "Fecha" (D) column
# so the newly-added rows (43-45) keep the same date formatting.
$dateFormat = $ws.Cells.Item(2,4).NumberFormat()

# Row 27
$ws.Cells.Item(27,1).Value = 9
$ws.Cells.Item(27,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(27,3).Value = 'Metropolitana'
$ws.Cells.Item(27,4).Value = 44467
$ws.Cells.Item(27,5).Value = 13
$ws.Cells.Item(27,6).Value = 300000000
$ws.Cells.Item(27,7).Value = 'Espárragos'
$ws.Cells.Item(27,8).Value = 'Sin especificar'
$ws.Cells.Item(27,9).Value = 'Banquete'
$ws.Cells.Item(27,10).Value = 106
$ws.Cells.Item(27,11).Value = 1700
$ws.Cells.Item(27,12).Value = 1800
$ws.Cells.Item(27,13).Value = 1750
$ws.Cells.Item(27,14).Value = '$/kilo'
$ws.Cells.Item(27,15).Value = 'Provincia de Linares'
$ws.Cells.Item(27,16).Value = 1750
$ws.Cells.Item(27,17).Value = 1
$ws.Cells.Item(27,18).Value = 'Hortaliza'

# Row 28
$ws.Cells.Item(28,1).Value = 9
$ws.Cells.Item(28,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(28,3).Value = 'Metropolitana'
$ws.Cells.Item(28,4).Value = 44467
$ws.Cells.Item(28,5).Value = 13
$ws.Cells.Item(28,6).Value = 300000000
$ws.Cells.Item(28,7).Value = 'Espárragos'
$ws.Cells.Item(28,8).Value = 'Sin especificar'
$ws.Cells.Item(28,9).Value = 'Primera'
$ws.Cells.Item(28,10).Value = 340
$ws.Cells.Item(28,11).Value = 1400
$ws.Cells.Item(28,12).Value = 1500
$ws.Cells.Item(28,13).Value = 1450
$ws.Cells.Item(28,14).Value = '$/kilo'
$ws.Cells.Item(28,15).Value = 'Provincia de Linares'
$ws.Cells.Item(28,16).Value = 1450
$ws.Cells.Item(28,17).Value = 1
$ws.Cells.Item(28,18).Value = 'Hortaliza'

# Row 29
$ws.Cells.Item(29,1).Value = 9
$ws.Cells.Item(29,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(29,3).Value = 'Metropolitana'
$ws.Cells.Item(29,4).Value = 44467
$ws.Cells.Item(29,5).Value = 13
$ws.Cells.Item(29,6).Value = 300000000
$ws.Cells.Item(29,7).Value = 'Espárragos'
$ws.Cells.Item(29,8).Value = 'Sin especificar'
$ws.Cells.Item(29,9).Value = 'Segunda'
$ws.Cells.Item(29,10).Value = 250
$ws.Cells.Item(29,11).Value = 1200
$ws.Cells.Item(29,12).Value = 1200
$ws.Cells.Item(29,13).Value = 1200
$ws.Cells.Item(29,14).Value = '$/kilo'
$ws.Cells.Item(29,15).Value = 'Provincia de Linares'
$ws.Cells.Item(29,16).Value = 1200
$ws.Cells.Item(29,17).Value = 1
$ws.Cells.Item(29,18).Value = 'Hortaliza'

# Row 30
$ws.Cells.Item(30,1).Value = 9
$ws.Cells.Item(30,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(30,3).Value = 'Metropolitana'
$ws.Cells.Item(30,4).Value = 44341
$ws.Cells.Item(30,5).Value = 13
$ws.Cells.Item(30,6).Value = 300000000
$ws.Cells.Item(30,7).Value = 'Espárragos'
$ws.Cells.Item(30,8).Value = 'Sin especificar'
$ws.Cells.Item(30,9).Value = 'Segunda'
$ws.Cells.Item(30,10).Value = 24
$ws.Cells.Item(30,11).Value = 28000
$ws.Cells.Item(30,12).Value = 30000
$ws.Cells.Item(30,13).Value = 29000
$ws.Cells.Item(30,14).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(30,15).Value = 'Región Metropolitana'
$ws.Cells.Item(30,16).Value = 2900
$ws.Cells.Item(30,17).Value = 10
$ws.Cells.Item(30,18).Value = 'Hortaliza'

# Row 31
$ws.Cells.Item(31,1).Value = 9
$ws.Cells.Item(31,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(31,3).Value = 'Metropolitana'
$ws.Cells.Item(31,4).Value = 44341
$ws.Cells.Item(31,5).Value = 13
$ws.Cells.Item(31,6).Value = 300000000
$ws.Cells.Item(31,7).Value = 'Espárragos'
$ws.Cells.Item(31,8).Value = 'Sin especificar'
$ws.Cells.Item(31,9).Value = 'Tercera'
$ws.Cells.Item(31,10).Value = 15
$ws.Cells.Item(31,11).Value = 24000
$ws.Cells.Item(31,12).Value = 26000
$ws.Cells.Item(31,13).Value = 25067
$ws.Cells.Item(31,14).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(31,15).Value = 'Región Metropolitana'
$ws.Cells.Item(31,16).Value = 2507
$ws.Cells.Item(31,17).Value = 10
$ws.Cells.Item(31,18).Value = 'Hortaliza'

# Row 32
$ws.Cells.Item(32,1).Value = 9
$ws.Cells.Item(32,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(32,3).Value = 'Metropolitana'
$ws.Cells.Item(32,4).Value = 44460
$ws.Cells.Item(32,5).Value = 13
$ws.Cells.Item(32,6).Value = 300000000
$ws.Cells.Item(32,7).Value = 'Espárragos'
$ws.Cells.Item(32,8).Value = 'Sin especificar'
$ws.Cells.Item(32,9).Value = 'Primera'
$ws.Cells.Item(32,10).Value = 250
$ws.Cells.Item(32,11).Value = 1900
$ws.Cells.Item(32,12).Value = 1900
$ws.Cells.Item(32,13).Value = 1900
$ws.Cells.Item(32,14).Value = '$/kilo'
$ws.Cells.Item(32,15).Value = 'Región Metropolitana'
$ws.Cells.Item(32,16).Value = 1900
$ws.Cells.Item(32,17).Value = 1
$ws.Cells.Item(32,18).Value = 'Hortaliza'

# Row 33
$ws.Cells.Item(33,1).Value = 9
$ws.Cells.Item(33,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(33,3).Value = 'Metropolitana'
$ws.Cells.Item(33,4).Value = 44460
$ws.Cells.Item(33,5).Value = 13
$ws.Cells.Item(33,6).Value = 300000000
$ws.Cells.Item(33,7).Value = 'Espárragos'
$ws.Cells.Item(33,8).Value = 'Sin especificar'
$ws.Cells.Item(33,9).Value = 'Segunda'
$ws.Cells.Item(33,10).Value = 160
$ws.Cells.Item(33,11).Value = 1700
$ws.Cells.Item(33,12).Value = 1700
$ws.Cells.Item(33,13).Value = 1700
$ws.Cells.Item(33,14).Value = '$/kilo'
$ws.Cells.Item(33,15).Value = 'Región Metropolitana'
$ws.Cells.Item(33,16).Value = 1700
$ws.Cells.Item(33,17).Value = 1
$ws.Cells.Item(33,18).Value = 'Hortaliza'

# Row 34
$ws.Cells.Item(34,1).Value = 9
$ws.Cells.Item(34,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(34,3).Value = 'Metropolitana'
$ws.Cells.Item(34,4).Value = 44460
$ws.Cells.Item(34,5).Value = 13
$ws.Cells.Item(34,6).Value = 300000000
$ws.Cells.Item(34,7).Value = 'Espárragos'
$ws.Cells.Item(34,8).Value = 'Sin especificar'
$ws.Cells.Item(34,9).Value = 'Tercera'
$ws.Cells.Item(34,10).Value = 106
$ws.Cells.Item(34,11).Value = 1500
$ws.Cells.Item(34,12).Value = 1500
$ws.Cells.Item(34,13).Value = 1500
$ws.Cells.Item(34,14).Value = '$/kilo'
$ws.Cells.Item(34,15).Value = 'Región Metropolitana'
$ws.Cells.Item(34,16).Value = 1500
$ws.Cells.Item(34,17).Value = 1
$ws.Cells.Item(34,18).Value = 'Hortaliza'

# Row 35
$ws.Cells.Item(35,1).Value = 9
$ws.Cells.Item(35,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(35,3).Value = 'Metropolitana'
$ws.Cells.Item(35,4).Value = 44466
$ws.Cells.Item(35,5).Value = 13
$ws.Cells.Item(35,6).Value = 300000000
$ws.Cells.Item(35,7).Value = 'Espárragos'
$ws.Cells.Item(35,8).Value = 'Sin especificar'
$ws.Cells.Item(35,9).Value = 'Banquete'
$ws.Cells.Item(35,10).Value = 110
$ws.Cells.Item(35,11).Value = 2500
$ws.Cells.Item(35,12).Value = 2500
$ws.Cells.Item(35,13).Value = 2500
$ws.Cells.Item(35,14).Value = '$/kilo'
$ws.Cells.Item(35,15).Value = 'Región Metropolitana'
$ws.Cells.Item(35,16).Value = 2500
$ws.Cells.Item(35,17).Value = 1
$ws.Cells.Item(35,18).Value = 'Hortaliza'

# Row 36
$ws.Cells.Item(36,1).Value = 9
$ws.Cells.Item(36,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(36,3).Value = 'Metropolitana'
$ws.Cells.Item(36,4).Value = 44463
$ws.Cells.Item(36,5).Value = 13
$ws.Cells.Item(36,6).Value = 300000000
$ws.Cells.Item(36,7).Value = 'Espárragos'
$ws.Cells.Item(36,8).Value = 'Sin especificar'
$ws.Cells.Item(36,9).Value = 'Banquete'
$ws.Cells.Item(36,10).Value = 34
$ws.Cells.Item(36,11).Value = 23000
$ws.Cells.Item(36,12).Value = 23000
$ws.Cells.Item(36,13).Value = 23000
$ws.Cells.Item(36,14).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(36,15).Value = 'Provincia de Linares'
$ws.Cells.Item(36,16).Value = 2300
$ws.Cells.Item(36,17).Value = 10
$ws.Cells.Item(36,18).Value = 'Hortaliza'

# Row 37
$ws.Cells.Item(37,1).Value = 9
$ws.Cells.Item(37,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(37,3).Value = 'Metropolitana'
$ws.Cells.Item(37,4).Value = 44463
$ws.Cells.Item(37,5).Value = 13
$ws.Cells.Item(37,6).Value = 300000000
$ws.Cells.Item(37,7).Value = 'Espárragos'
$ws.Cells.Item(37,8).Value = 'Sin especificar'
$ws.Cells.Item(37,9).Value = 'Primera'
$ws.Cells.Item(37,10).Value = 52
$ws.Cells.Item(37,11).Value = 21000
$ws.Cells.Item(37,12).Value = 21000
$ws.Cells.Item(37,13).Value = 21000
$ws.Cells.Item(37,14).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(37,15).Value = 'Provincia de Linares'
$ws.Cells.Item(37,16).Value = 2100
$ws.Cells.Item(37,17).Value = 10
$ws.Cells.Item(37,18).Value = 'Hortaliza'

# Row 38
$ws.Cells.Item(38,1).Value = 9
$ws.Cells.Item(38,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(38,3).Value = 'Metropolitana'
$ws.Cells.Item(38,4).Value = 44463
$ws.Cells.Item(38,5).Value = 13
$ws.Cells.Item(38,6).Value = 300000000
$ws.Cells.Item(38,7).Value = 'Espárragos'
$ws.Cells.Item(38,8).Value = 'Sin especificar'
$ws.Cells.Item(38,9).Value = 'Primera'
$ws.Cells.Item(38,10).Value = 340
$ws.Cells.Item(38,11).Value = 2000
$ws.Cells.Item(38,12).Value = 2000
$ws.Cells.Item(38,13).Value = 2000
$ws.Cells.Item(38,14).Value = '$/kilo'
$ws.Cells.Item(38,15).Value = 'Provincia de Linares'
$ws.Cells.Item(38,16).Value = 2000
$ws.Cells.Item(38,17).Value = 1
$ws.Cells.Item(38,18).Value = 'Hortaliza'

# Row 39
$ws.Cells.Item(39,1).Value = 9
$ws.Cells.Item(39,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(39,3).Value = 'Metropolitana'
$ws.Cells.Item(39,4).Value = 44463
$ws.Cells.Item(39,5).Value = 13
$ws.Cells.Item(39,6).Value = 300000000
$ws.Cells.Item(39,7).Value = 'Espárragos'
$ws.Cells.Item(39,8).Value = 'Sin especificar'
$ws.Cells.Item(39,9).Value = 'Segunda'
$ws.Cells.Item(39,10).Value = 43
$ws.Cells.Item(39,11).Value = 18000
$ws.Cells.Item(39,12).Value = 18000
$ws.Cells.Item(39,13).Value = 18000
$ws.Cells.Item(39,14).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(39,15).Value = 'Provincia de Linares'
$ws.Cells.Item(39,16).Value = 1800
$ws.Cells.Item(39,17).Value = 10
$ws.Cells.Item(39,18).Value = 'Hortaliza'

# Row 40
$ws.Cells.Item(40,1).Value = 9
$ws.Cells.Item(40,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(40,3).Value = 'Metropolitana'
$ws.Cells.Item(40,4).Value = 44463
$ws.Cells.Item(40,5).Value = 13
$ws.Cells.Item(40,6).Value = 300000000
$ws.Cells.Item(40,7).Value = 'Espárragos'
$ws.Cells.Item(40,8).Value = 'Sin especificar'
$ws.Cells.Item(40,9).Value = 'Segunda'
$ws.Cells.Item(40,10).Value = 160
$ws.Cells.Item(40,11).Value = 1700
$ws.Cells.Item(40,12).Value = 1700
$ws.Cells.Item(40,13).Value = 1700
$ws.Cells.Item(40,14).Value = '$/kilo'
$ws.Cells.Item(40,15).Value = 'Provincia de Linares'
$ws.Cells.Item(40,16).Value = 1700
$ws.Cells.Item(40,17).Value = 1
$ws.Cells.Item(40,18).Value = 'Hortaliza'

# Row 41
$ws.Cells.Item(41,1).Value = 9
$ws.Cells.Item(41,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(41,3).Value = 'Metropolitana'
$ws.Cells.Item(41,4).Value = 44160
$ws.Cells.Item(41,5).Value = 13
$ws.Cells.Item(41,6).Value = 300000000
$ws.Cells.Item(41,7).Value = 'Espárragos'
$ws.Cells.Item(41,8).Value = 'Verde'
$ws.Cells.Item(41,9).Value = 'Banquete'
$ws.Cells.Item(41,10).Value = 150
$ws.Cells.Item(41,11).Value = 1000
$ws.Cells.Item(41,12).Value = 1050
$ws.Cells.Item(41,13).Value = 1033
$ws.Cells.Item(41,14).Value = '$/kilo'
$ws.Cells.Item(41,15).Value = 'Región Metropolitana'
$ws.Cells.Item(41,16).Value = 1033
$ws.Cells.Item(41,17).Value = 1
$ws.Cells.Item(41,18).Value = 'Hortaliza'

# Row 42
$ws.Cells.Item(42,1).Value = 9
$ws.Cells.Item(42,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(42,3).Value = 'Metropolitana'
$ws.Cells.Item(42,4).Value = 44160
$ws.Cells.Item(42,5).Value = 13
$ws.Cells.Item(42,6).Value = 300000000
$ws.Cells.Item(42,7).Value = 'Espárragos'
$ws.Cells.Item(42,8).Value = 'Verde'
$ws.Cells.Item(42,9).Value = 'Primera'
$ws.Cells.Item(42,10).Value = 220
$ws.Cells.Item(42,11).Value = 800
$ws.Cells.Item(42,12).Value = 800
$ws.Cells.Item(42,13).Value = 800
$ws.Cells.Item(42,14).Value = '$/kilo'
$ws.Cells.Item(42,15).Value = 'Región Metropolitana'
$ws.Cells.Item(42,16).Value = 800
$ws.Cells.Item(42,17).Value = 1
$ws.Cells.Item(42,18).Value = 'Hortaliza'

# Row 43
$ws.Cells.Item(43,1).Value = 9
$ws.Cells.Item(43,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(43,3).Value = 'Metropolitana'
$ws.Cells.Item(43,4).Value = 44160
$ws.Cells.Item(43,4).NumberFormat = $dateFormat
$ws.Cells.Item(43,5).Value = 13
$ws.Cells.Item(43,6).Value = 300000000
$ws.Cells.Item(43,7).Value = 'Espárragos'
$ws.Cells.Item(43,8).Value = 'Verde'
$ws.Cells.Item(43,9).Value = 'Segunda'
$ws.Cells.Item(43,10).Value = 160
$ws.Cells.Item(43,11).Value = 650
$ws.Cells.Item(43,12).Value = 650
$ws.Cells.Item(43,13).Value = 650
$ws.Cells.Item(43,14).Value = '$/kilo'
$ws.Cells.Item(43,15).Value = 'Región Metropolitana'
$ws.Cells.Item(43,16).Value = 650
$ws.Cells.Item(43,17).Value = 1
$ws.Cells.Item(43,18).Value = 'Hortaliza'

# Row 44
$ws.Cells.Item(44,1).Value = 9
$ws.Cells.Item(44,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(44,3).Value = 'Metropolitana'
$ws.Cells.Item(44,4).Value = 44176
$ws.Cells.Item(44,4).NumberFormat = $dateFormat
$ws.Cells.Item(44,5).Value = 13
$ws.Cells.Item(44,6).Value = 300000000
$ws.Cells.Item(44,7).Value = 'Espárragos'
$ws.Cells.Item(44,8).Value = 'Sin especificar'
$ws.Cells.Item(44,9).Value = 'Primera'
$ws.Cells.Item(44,10).Value = 4300
$ws.Cells.Item(44,11).Value = 1300
$ws.Cells.Item(44,12).Value = 1300
$ws.Cells.Item(44,13).Value = 1300
$ws.Cells.Item(44,14).Value = '$/kilo'
$ws.Cells.Item(44,15).Value = 'Provincia de Linares'
$ws.Cells.Item(44,16).Value = 1300
$ws.Cells.Item(44,17).Value = 1
$ws.Cells.Item(44,18).Value = 'Hortaliza'

# Row 45
$ws.Cells.Item(45,1).Value = 9
$ws.Cells.Item(45,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(45,3).Value = 'Metropolitana'
$ws.Cells.Item(45,4).Value = 44176
$ws.Cells.Item(45,4).NumberFormat = $dateFormat
$ws.Cells.Item(45,5).Value = 13
$ws.Cells.Item(45,6).Value = 300000000
$ws.Cells.Item(45,7).Value = 'Espárragos'
$ws.Cells.Item(45,8).Value = 'Sin especificar'
$ws.Cells.Item(45,9).Value = 'Segunda'
$ws.Cells.Item(45,10).Value = 2500
$ws.Cells.Item(45,11).Value = 1000
$ws.Cells.Item(45,12).Value = 1000
$ws.Cells.Item(45,13).Value = 1000
$ws.Cells.Item(45,14).Value = '$/kilo'
$ws.Cells.Item(45,15).Value = 'Provincia de Linares'
$ws.Cells.Item(45,16).Value = 1000
$ws.Cells.Item(45,17).Value = 1
$ws.Cells.Item(45,18).Value = 'Hortaliza'
